$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 99

$ws.Cells.Item($newRow, 1).Value = "Need: Golang System integration Specialist: Denver, CO OR West Chester, PA (Onsite from Day 1)"
$ws.Cells.Item($newRow, 2).Value = "https://www.dice.com/job-detail/b698b9ba-5c0e-4766-aee6-47ddeedd9fd7"
$ws.Cells.Item($newRow, 3).Value = "Denver, Colorado"
$ws.Cells.Item($newRow, 4).Value = "Third Party, Contract"
$ws.Cells.Item($newRow, 5).Value = "Depends on Experience"
$ws.Cells.Item($newRow, 6).Value = "Shrive Technologies LLC"
